# Auto-generated PowerShell script to apply crypto price/volume updates
# Commit message: Updated cryptos list on Tue Feb 20 13:47:28 UTC 2024 with GitHub Actions
#
# Source data is rendered as plain text cells (t="inlineStr") in the workbook.
# Some of the replacement strings (e.g. "1.00", "2.80") would otherwise be
# silently re-interpreted as numbers by Excel's Range.Value setter, dropping
# the trailing zero / the original text formatting. For those specific cells
# we prefix the literal with an apostrophe (Excel's "treat as text" marker)
# and then restore the cell's style to "Normal" so no stray number-format
# style is left attached to the cell - matching the source which has no
# explicit style on these cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.859.89'
$ws.Range('E2').Value = '  +1.23%  '
$ws.Range('D3').Value = '2.992.71'
$ws.Range('E3').Value = '  +3.15%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'359.45"
$ws.Range('E5').Value = '  +1.90%  '
$ws.Range('D6').Value = "'110.68"
$ws.Range('E6').Value = '  -2.21%  '
$ws.Range('E7').Value = '  +3.22%  '
$ws.Range('D8').Value = "'1.00"
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = "'0.638"
$ws.Range('E9').Value = '  +2.19%  '
$ws.Range('D10').Value = "'39.54"
$ws.Range('E10').Value = '  -1.22%  '
$ws.Range('D11').Value = "'0.0884"
$ws.Range('E11').Value = '  +2.72%  '
$ws.Range('E12').Value = '  +1.69%  '
$ws.Range('D13').Value = "'19.66"
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'7.86"
$ws.Range('E14').Value = '  +1.48%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.441.18'
$ws.Range('E15').Value = '  +2.48%  '
$ws.Range('D16').Value = '2.975.28'
$ws.Range('E16').Value = '  +2.06%  '
$ws.Range('D17').Value = "'0.997"
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('D18').Value = '52.848.77'
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').Value = "'3.51"
$ws.Range('E19').Value = '  +6.19%  '
$ws.Range('E20').Value = '  +0.68%  '
$ws.Range('D21').Value = "'14.03"
$ws.Range('E21').Value = '  -1.11%  '
$ws.Range('D22').Value = '0.0₃0991'
$ws.Range('E22').Value = '  +1.51%  '
$ws.Range('D23').Value = "'273.38"
$ws.Range('E23').Value = '  +1.50%  '
$ws.Range('D24').Value = "'70.76"
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('D25').Value = "'2.83"
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('D26').Value = "'0.181"
$ws.Range('E26').Value = '  +3.86%  '
$ws.Range('D27').Value = "'7.86"
$ws.Range('E27').Value = '  +18.36%  '
$ws.Range('D28').Value = "'27.49"
$ws.Range('E28').Value = '  +2.71%  '
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('D30').Value = "'0.108"
$ws.Range('E30').Value = '  +5.78%  '
$ws.Range('D31').Value = "'10.58"
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('D32').Value = "'38.42"
$ws.Range('E32').Value = '  +2.34%  '
$ws.Range('D33').Value = "'6.17"
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('E34').Value = '  +10.67%  '
$ws.Range('D35').Value = "'52.46"
$ws.Range('E35').Value = '  -1.42%  '
$ws.Range('D36').Value = "'0.0448"
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('D37').Value = "'0.996"
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').Value = "'3.28"
$ws.Range('E38').Value = '  -1.34%  '
$ws.Range('D39').Value = "'2.05"
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('D40').Value = "'18.37"
$ws.Range('E40').Value = '  -2.60%  '
$ws.Range('D41').Value = "'2.80"
$ws.Range('E41').Value = '  +1.66%  '
$ws.Range('E42').Value = '  +3.03%  '
$ws.Range('D43').Value = "'23.91"
$ws.Range('E43').Value = '  +3.91%  '
$ws.Range('D44').Value = "'119.89"
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('D45').Value = "'2.17"
$ws.Range('E45').Value = '  -1.13%  '
$ws.Range('D46').Value = "'3.51"
$ws.Range('E46').Value = '  -0.19%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.156.68'
$ws.Range('E47').Value = '  -0.98%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = "'2.47"
$ws.Range('E48').Value = '  -5.05%  '
$ws.Range('D49').Value = "'0.0358"
$ws.Range('E49').Value = '  +2.53%  '
$ws.Range('E50').Value = '  -5.04%  '
$ws.Range('D51').Value = "'0.926"
$ws.Range('E51').Value = '  -2.64%  '

# Restore default ('Normal') style on cells that needed the text-marker
# apostrophe, so no residual number-format style is attached to them.
$ws.Range('D4').Style = "Normal"
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D13').Style = "Normal"
$ws.Range('D14').Style = "Normal"
$ws.Range('D17').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D51').Style = "Normal"
